$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$f = $ws.Range("C10:G10").Font
$f.Name = "Bookman Old Style"
$f.Size = 13
$f.Bold = $true
$f.Underline = $true
$ws.Range("C10:G10").HorizontalAlignment = -4108
